# Update cryptos list data (prices and 1h volume-change percentages)
# Mirrors a scheduled data refresh (GitHub Actions cron job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.516.18"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "3.099.78"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.90"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.091.95"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  +6.70%  "
$ws.Range("E11").Value = "  -2.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.456"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.31"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.81%  "
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "3.612.57"
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "63.348.31"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.09"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").Value = "3.091.35"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.91"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.23"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.725"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.46"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("E24").Value = "  -2.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.31"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.03"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +9.19%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.67"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.84"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.63"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").Value = "0.0₃0853"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.39"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.45%  "
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.03"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.99"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.21"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "435.12"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0368"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "2.875.18"
$ws.Range("E44").Value = "  -1.92%  "
$ws.Range("E45").Value = "  -3.15%  "
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.03"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.57"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.110"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.12"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.13%  "
